$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Cells.Item(3, 4).Style

# Row 2
$ws.Cells.Item(2, 4).Value = '60.245.74'
$ws.Cells.Item(2, 5).Value = '  -2.72%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.384.63'
$ws.Cells.Item(3, 5).Value = '  -4.52%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '538.67'
$ws.Cells.Item(5, 5).Value = '  -2.22%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '139.20'
$ws.Cells.Item(6, 4).Style = $plainStyle
$ws.Cells.Item(6, 5).Value = '  -5.23%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.06%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.574'
$ws.Cells.Item(8, 5).Value = '  -6.84%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.385.86'
$ws.Cells.Item(9, 5).Value = '  -4.53%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.104'
$ws.Cells.Item(10, 5).Value = '  -3.22%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.08%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '5.31'
$ws.Cells.Item(12, 5).Value = '  -1.41%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '0.339'
$ws.Cells.Item(13, 5).Value = '  -4.75%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '25.17'
$ws.Cells.Item(14, 5).Value = '  -3.73%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.819.13'
$ws.Cells.Item(15, 5).Value = '  -4.43%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000163'
$ws.Cells.Item(16, 4).Style = $plainStyle
$ws.Cells.Item(16, 5).Value = '  +0.20%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '60.239.09'
$ws.Cells.Item(17, 5).Value = '  -2.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.388.52'
$ws.Cells.Item(18, 5).Value = '  -4.44%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '10.56'
$ws.Cells.Item(19, 5).Value = '  -4.82%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '6.67'
$ws.Cells.Item(20, 5).Value = '  -5.08%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(21, 4).Value = '4.04'
$ws.Cells.Item(21, 5).Value = '  -3.58%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '312.52'
$ws.Cells.Item(22, 5).Value = '  -2.76%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.03%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '1.78'
$ws.Cells.Item(24, 5).Value = '  +2.68%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '62.57'
$ws.Cells.Item(25, 5).Value = '  -2.10%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 4).Style = $plainStyle
$ws.Cells.Item(26, 5).Value = '  +0.49%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '2.503.72'
$ws.Cells.Item(27, 5).Value = '  -4.55%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Aptos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.60'
$ws.Cells.Item(28, 4).Style = $plainStyle
$ws.Cells.Item(28, 5).Value = '  -0.33%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).Value = '0.0₃0901'
$ws.Cells.Item(29, 5).Value = '  -10.70%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '7.98'
$ws.Cells.Item(30, 5).Value = '  -4.81%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '1.41'
$ws.Cells.Item(31, 5).Value = '  -5.25%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '501.47'
$ws.Cells.Item(32, 5).Value = '  -6.13%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -4.48%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.82'
$ws.Cells.Item(34, 5).Value = '  -3.98%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.55'
$ws.Cells.Item(35, 5).Value = '  -1.12%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.07%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.60'
$ws.Cells.Item(37, 4).Style = $plainStyle
$ws.Cells.Item(37, 5).Value = '  -5.57%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '5.40'
$ws.Cells.Item(38, 4).Style = $plainStyle
$ws.Cells.Item(38, 5).Value = '  -8.03%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.371'
$ws.Cells.Item(39, 5).Value = '  -1.89%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '17.90'
$ws.Cells.Item(40, 4).Style = $plainStyle
$ws.Cells.Item(40, 5).Value = '  -3.25%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.10%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '137.11'
$ws.Cells.Item(42, 5).Value = '  -4.74%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '1.69'
$ws.Cells.Item(43, 5).Value = '  -0.26%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '40.30'
$ws.Cells.Item(44, 4).Style = $plainStyle
$ws.Cells.Item(44, 5).Value = '  -0.05%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.15'
$ws.Cells.Item(45, 5).Value = '  -6.28%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '139.43'
$ws.Cells.Item(46, 5).Value = '  -6.55%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.50'
$ws.Cells.Item(47, 4).Style = $plainStyle
$ws.Cells.Item(47, 5).Value = '  -2.15%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '20.05'
$ws.Cells.Item(48, 5).Value = '  -3.58%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.0513'
$ws.Cells.Item(49, 5).Value = '  -4.21%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.575'
$ws.Cells.Item(50, 5).Value = '  -2.46%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.0922'
$ws.Cells.Item(51, 5).Value = '  -3.00%  '
